# Updates cryptocurrency price (D) and 1h volume-change (E) figures
# to match the refreshed coinranking.com scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.557.77'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.05%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.603.48'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.44%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.39'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.04%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.522'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.39%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '26.78'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +7.14%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '43.48'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.48%  '

$ws.Range("E10").Value = '  +2.12%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0601'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.45%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0907'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.30%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.833.23'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.47%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.609.32'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.88%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.559.13'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.91%  '

$ws.Range("E16").Value = '  +3.56%  '

$ws.Range("E17").Value = '  +2.11%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '63.42'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.07%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '240.40'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.47%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.59'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.31%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0691'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.18%  '

$ws.Range("E22").Value = '  +0.12%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.98'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.27%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.21'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.93%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.08'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.54%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.15'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.48%  '

$ws.Range("E27").Value = '  +2.68%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.25'
$ws.Range("D28").Style = "Normal"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.39'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.35%  '

$ws.Range("E30").Value = '  +0.10%  '

$ws.Range("E31").Value = '  +2.65%  '

$ws.Range("E32").Value = '  +0.60%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.22'
$ws.Range("D33").Style = "Normal"

$ws.Range("E34").Value = '  +3.59%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.407.17'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.42%  '

$ws.Range("E36").Value = '  +0.36%  '

$ws.Range("E37").Value = '  +4.69%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.83'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.35%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.31'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.31%  '

$ws.Range("E40").Value = '  +2.35%  '

$ws.Range("E41").Value = '  +3.68%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.99'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.32%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0494'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +7.96%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '53.80'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +26.57%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.796'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.47%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.03%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '65.75'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.92%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.28'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.06%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.744.76'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.16%  '

$ws.Range("E50").Value = '  -1.76%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '86.58'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.03%  '
